# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the crypto symbol table, per commit "Updated symbol list on Mon Jan 16
# 09:41:23 UTC 2023 with GitHub Actions".
#
# Each new value is written with a leading apostrophe so Excel stores it as
# literal text (matching the original cell type) instead of silently
# reinterpreting the numeric-looking / percentage-looking text as a number.
# The apostrophe is a type-marker consumed by Excel, not stored in the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''299.56'
$ws.Range("E2").Value = '''1.66%'
$ws.Range("D3").Value = '''31.21'
$ws.Range("E3").Value = '''-0.41%'
$ws.Range("D4").Value = '''5.130'
$ws.Range("E4").Value = '''0.41%'
$ws.Range("D5").Value = '''0.08111'
$ws.Range("E5").Value = '''10.19%'
$ws.Range("D6").Value = '''2.730'
$ws.Range("E6").Value = '''66.37%'
$ws.Range("D7").Value = '''7.850'
$ws.Range("E7").Value = '''2.32%'
$ws.Range("D8").Value = '''3.833'
$ws.Range("E8").Value = '''1.93%'
$ws.Range("D9").Value = '''0.9089'
$ws.Range("E9").Value = '''-1.24%'
$ws.Range("D10").Value = '''0.1719'
$ws.Range("E10").Value = '''2.72%'
$ws.Range("D11").Value = '''0.07259'
$ws.Range("E11").Value = '''2.90%'
$ws.Range("D12").Value = '''0.07970'
$ws.Range("E12").Value = '''-0.12%'
$ws.Range("D13").Value = '''0.03025'
$ws.Range("E13").Value = '''0.89%'
$ws.Range("D14").Value = '''0.09971'
$ws.Range("E14").Value = '''0.86%'
$ws.Range("D15").Value = '''0.001499'
$ws.Range("E15").Value = '''0.72%'
$ws.Range("D16").Value = '''0.006015'
$ws.Range("E16").Value = '''-2.90%'
$ws.Range("D17").Value = '''3.497'
$ws.Range("E17").Value = '''1.32%'
$ws.Range("D18").Value = '''2.254'
$ws.Range("E18").Value = '''1.14%'
$ws.Range("D20").Value = '''0.1339'
$ws.Range("E20").Value = '''0.35%'
$ws.Range("D21").Value = '''4.602'
$ws.Range("E21").Value = '''0.89%'
$ws.Range("E22").Value = '''3.31%'
$ws.Range("D23").Value = '''0.04571'
$ws.Range("E23").Value = '''-1.18%'
$ws.Range("D24").Value = '''0.001264'
$ws.Range("D25").Value = '''0.004446'
$ws.Range("E25").Value = '''0.64%'
$ws.Range("E26").Value = '''-9.17%'
$ws.Range("D27").Value = '''0.0003434'
$ws.Range("E27").Value = '''83.01%'
$ws.Range("D39").Value = '''0.01817'
$ws.Range("E39").Value = '''8.32%'
$ws.Range("D40").Value = '''0.04537'
$ws.Range("E40").Value = '''3.06%'
$ws.Range("D41").Value = '''0.007053'
$ws.Range("E41").Value = '''-0.85%'
$ws.Range("E42").Value = '''1.19%'
$ws.Range("D43").Value = '''0.002243'
$ws.Range("E43").Value = '''6.22%'
$ws.Range("D44").Value = '''0.01051'
$ws.Range("E44").Value = '''-4.34%'
$ws.Range("D45").Value = '''0.00006310'
$ws.Range("E45").Value = '''5.58%'
$ws.Range("E46").Value = '''0.09%'
$ws.Range("D47").Value = '''0.006407'
$ws.Range("E47").Value = '''-41.77%'
$ws.Range("E48").Value = '''15.31%'
$ws.Range("D49").Value = '''0.00002102'
$ws.Range("E49").Value = '''0.09%'
$ws.Range("D50").Value = '''0.0002002'
$ws.Range("E50").Value = '''0.16%'
